$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates (card holder name / card number)
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay TEXT (not be coerced into a
# Number, which would both change the cell's type and lose precision on a
# 16-digit value since doubles only carry ~15 significant digits). Build
# the text in an unused scratch cell via a formula that forces a string
# result, copy it, and paste VALUES ONLY into B3 so the existing number
# formatting/style of B3 (s="8") is left completely untouched.
$scratch = $ws.Range("H1")
$scratch.Formula = "=""2570314725427075"""
$scratch.Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$scratch.Clear()

$ws.Range("C3").Value = "Mohaupt"

# Opening balance label
$ws.Range("D5").Value = "KONTOSTAND AM 09.07.2025"

# Row 6
$ws.Range("B6").Value = "12.07."
$ws.Range("C6").Value = "13.07."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 3954591"
$ws.Range("E6").Value = "83,87-"

# Row 7
$ws.Range("B7").Value = "16.07."
$ws.Range("C7").Value = "17.07."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 33900730"
$ws.Range("E7").Value = "39,76-"

# Row 8
$ws.Range("B8").Value = "20.07."
$ws.Range("C8").Value = "21.07."
$ws.Range("D8").Value = "KARTENZ./20.07 REWE RO"
$ws.Range("E8").Value = "95,32-"

# Rows 9-11 no longer hold transactions; blank them out (values removed,
# alignment on the amount column differs per row as in the template)
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Closing balance / next statement date
$ws.Range("D12").Value = "KONTOSTAND AM 25.07.2025"
$ws.Range("E12").Value = "218,95-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.07.2025"
